# The commit swaps the OOXML content of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: theme1.xml (the Slide Master's theme, currently
# the "Integral" colour scheme) ends up holding the stock "Office Theme"
# colour scheme, while theme2.xml (the Notes Master's theme) ends up
# holding the "Integral" colours that used to live in theme1.xml. The
# font scheme and format scheme (fills/lines/effects) are identical
# between the two themes in this deck, so the only observable change is
# the 12-colour theme colour scheme used by the Slide Master.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

# RGB() below packs as 0x00BBGGRR, matching the PowerPoint COM
# convention used by RGBColor.RGB / ColorFormat.RGB.
function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target values: the stock PowerPoint "Office Theme" colour scheme.
# Index order (MsoThemeColorSchemeIndex): 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
# 11 hyperlink, 12 followed hyperlink.
$themeColors.Item(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1      000000
$themeColors.Item(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1      FFFFFF
$themeColors.Item(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2      44546A
$themeColors.Item(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2      E7E6E6
$themeColors.Item(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1  5B9BD5
$themeColors.Item(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2  ED7D31
$themeColors.Item(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3  A5A5A5
$themeColors.Item(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4  FFC000
$themeColors.Item(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5  4472C4
$themeColors.Item(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6  70AD47
$themeColors.Item(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink    0563C1
$themeColors.Item(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink 954F72
